$wb = $excel.ActiveWorkbook

# Component placement (CPL) tweak for Q2 on the DuDad-bottom-pos sheet:
# nudge "Mid X" from 54.2068 to 54.3068 (part of the resistor-value /
# LED-dimming rework referenced in the commit).
$ws = $wb.Worksheets.Item("DuDad-bottom-pos")
$ws.Range("B5").Value = 54.306800000000003

# Drop the stray, empty "Sheet1" tab left over from the BOM-generation work.
$sheet1 = $wb.Worksheets.Item("Sheet1")
[void]$sheet1.Delete()
